$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3521, 4201, 4272, 4331, 4331, 4331, 4374, 4374, 4407, 4407, 4407, 4559, 4641, 4641)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
